# Calendario-EP-2026-2: add "Asistencia" sheet (attendance tracker) after "Hoja1"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet, positioned right after Hoja1 -------------------
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Asistencia"

# --- Pull the bold/gray header style from Hoja1!A1:C1 (s=1) ---------------
$ws1.Range("A1:C1").Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)

# --- Header row: Nombre + 9 fortnightly dates (Tue/Fri pairs) -------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = 46056
$ws.Range("C1").Value = 46059
$ws.Range("D1").Formula = "=B1+7"
$ws.Range("E1").Formula = "=C1+7"
$ws.Range("F1").Formula = "=D1+7"
$ws.Range("G1").Formula = "=E1+7"
$ws.Range("H1").Formula = "=F1+7"
$ws.Range("I1").Formula = "=G1+7"
$ws.Range("J1").Formula = "=H1+7"
$ws.Range("B1:J1").NumberFormat = "dd/mm/yy"

# --- Student names + attendance values -------------------------------------
$ws.Range("A2").Value = "Ana Cristina Uc Canela"
$ws.Range("A3").Value = "Xiadani Briones García"
$ws.Range("A4").Value = "Jamil Jassiel Hernández Enríquez"
$ws.Range("A5").Value = "Casandra Gallardo Badillo"
$ws.Range("A6").Value = "Renata Flores García"
$ws.Range("A7").Value = "Pablo Villanueva Valdez"
$ws.Range("A8").Value = "Rodrigo Chan Catzim"

$ws.Range("B2:H2").Value = 1
$ws.Range("B3:H3").Value = 1
$ws.Range("B4:H4").Value = 1
$ws.Range("B5:H5").Value = 1

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1

# --- Column width + zoom on the new sheet ----------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.72

# --- View state: both sheets zoomed to 180%, Asistencia becomes active ----
$ws1.Activate()
$excel.ActiveWindow.Zoom = 180

$ws.Activate()
$excel.ActiveWindow.Zoom = 180
$ws.Range("H3").Select()
